$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph by its index,
# since newly-inserted paragraphs are most reliably addressed via
# $d.Paragraphs.Item(n) rather than chained .Next references.
$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $targetIdx = $idx
    }
}

if ($targetIdx -gt 0) {
    $docentePara = $d.Paragraphs.Item($targetIdx)
    $docentePara.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIdx + 1)
    $newPara.Range.Text = "5840942 - Marco Aurélio Kondracki de Alcântara"
    $newPara.Style = "ListBullet"
}
